# Update weekly Fruta/Hortaliza price rows: reorder the data among rows
# 3,4,5,6,8,9,10 (row 1 header, row 2 and row 7 unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for each affected row: Date (D), Volumen (M), Precio minimo (N),
# Precio maximo (O), Precio promedio ponderado (P), Precio $/Kg (S)
$updates = @{
    3  = @{ D = 44482; M = 240; N = 10000; O = 11000; P = 10500; S = 5250 }
    4  = @{ D = 44517; M = 400; N = 5500;  O = 6000;  P = 5750;  S = 2875 }
    5  = @{ D = 44461; M = 200; N = 11000; O = 12000; P = 11500; S = 5750 }
    6  = @{ D = 44490; M = 400; N = 9500;  O = 10000; P = 9750;  S = 4875 }
    8  = @{ D = 44489; M = 160; N = 9500;  O = 10000; P = 9750;  S = 4875 }
    9  = @{ D = 44497; M = 500; N = 9000;  O = 10000; P = 9500;  S = 4750 }
    10 = @{ D = 44455; M = 200; N = 12000; O = 13000; P = 12500; S = 6250 }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("N$row").Value = $vals.N
    $ws.Range("O$row").Value = $vals.O
    $ws.Range("P$row").Value = $vals.P
    $ws.Range("S$row").Value = $vals.S
}
